# Results.xlsx — "Completed 1.1 and 1.2"
# Adds benchmark rows for the "leaky" activation function and a "pca"
# (10 components) run, and highlights the best (min Avg) row for the
# sigmoid group (row 12) the same way the existing relu/tanh "best" rows
# (4 and 9) are already highlighted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# --- 1) Highlight sigmoid's best row (row 12, lr=0.1, min Avg) -------------
# Re-use the existing highlighted-row format (as seen on rows 4 and 9)
# instead of poking Interior.Color directly, so the workbook keeps reusing
# the same cell style rather than minting a duplicate one.
$ws.Range("A4:E4").Copy()
$ws.Range("A12:E12").PasteSpecial($xlPasteFormats)
$ws.Application.CutCopyMode = $false

# --- 2) New "leaky" rows (17-21) --------------------------------------------
$ws.Range("A17").Value = "leaky"
$ws.Range("B17").Value = 0.1
$ws.Range("C17").Value = 0.000282
$ws.Range("D17").Value = 0.050595
$ws.Range("E17").Value = 0.010952

$ws.Range("A18").Value = "leaky"
$ws.Range("B18").Value = 0.01
$ws.Range("C18").Value = 0.000223
$ws.Range("D18").Value = 0.00025
$ws.Range("E18").Value = 0.000233

$ws.Range("A19").Value = "leaky"
$ws.Range("B19").Value = 0.001
$ws.Range("C19").Value = 0.000189
$ws.Range("D19").Value = 0.000231
$ws.Range("E19").Value = 0.000224

$ws.Range("A20").Value = "leaky"
$ws.Range("B20").Value = 0.0001
$ws.Range("C20").Value = 0.000238
$ws.Range("D20").Value = 0.000241
$ws.Range("E20").Value = 0.00024

$ws.Range("A21").Value = "leaky"
$ws.Range("B21").Value = 0.00001
$ws.Range("C21").Value = 0.000268
$ws.Range("D21").Value = 0.000272
$ws.Range("E21").Value = 0.00027

# Highlight leaky's best row (row 19, lr=0.001, min Avg) the same way.
$ws.Range("A4:E4").Copy()
$ws.Range("A19:E19").PasteSpecial($xlPasteFormats)
$ws.Application.CutCopyMode = $false

# --- 3) New "pca" row (22) --------------------------------------------------
# Only A/B/E are populated (C/D — Min/Max — don't apply to a single PCA run).
$ws.Range("A22").Value = "pca"
$ws.Range("B22").Value = "10 components"
$ws.Range("E22").Value = 0.000142

# pca only has one run, so it is highlighted too (consistent with the
# "best row per function" highlighting convention used throughout). Paste
# the format cell-by-cell so empty C22/D22 aren't materialised.
$ws.Range("A4").Copy()
$ws.Range("A22").PasteSpecial($xlPasteFormats)
$ws.Range("B4").Copy()
$ws.Range("B22").PasteSpecial($xlPasteFormats)
$ws.Range("E4").Copy()
$ws.Range("E22").PasteSpecial($xlPasteFormats)
$ws.Application.CutCopyMode = $false

# --- 4) View state: scroll back to the top, select C22 ----------------------
$null = $ws.Range("C22").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# --- 5) Column widths (minor autofit-style narrowing seen in the commit) ---
$ws.Columns.Item(1).ColumnWidth = 10.508503401360565
$ws.Columns.Item(2).ColumnWidth = 13.477891156462565
